$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("thermal_hull")

# Rename "Außenwand (netto)" -> "Aussenwand"
$ws.Range("A2").Value = "Aussenwand"

# Update the active selection on the sheet to A3
$ws.Activate()
$ws.Range("A3").Select()
